# Updated Backplate Stud Size
# Changed stud size from 6-32 to 8-32
#
# Add the new ECO (Engineering Change Order) log entry row to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of change-log data (row 3)
$ws.Range("A3").Value = "POSCON.v1.1"
$ws.Range("B3").Value = "Fixed backplate stud size"
$ws.Range("C3").Value = "Marc Levinson"
$ws.Range("D3").Value = "Inconcistent labeling between drawing and BOM"
$ws.Range("E3").Value = 41927

# Match the Date column's existing date formatting (copy format from E2)
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Widen the "Change" column so the new, longer text is readable
$ws.Columns.Item(2).ColumnWidth = 18.1425

# Move the active selection, as recorded after the edit
$ws.Range("D9").Select()
